$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.134548783302307
$ws.Range("B1").Value = 2.158451080322266
$ws.Range("C1").Value = 2.793316125869751
$ws.Range("D1").Value = 1.482792377471924
$ws.Range("E1").Value = 0.9023613333702087
